$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells in column D hold numeric-looking text (t="inlineStr"), so the
# new values must be written back as text, not auto-converted to numbers.
# Trick: temporarily mark the cell as Text ("@") so Excel accepts the
# numeric-looking literal as a string, then restore the cell's original
# "Normal" style so no stray formatting is left behind.
function Set-TextValue {
    param($cell, $value)
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

# Price (column D) updates
Set-TextValue $ws.Range("D2")  "244.70"
Set-TextValue $ws.Range("D4")  "5.410"
Set-TextValue $ws.Range("D5")  "0.05990"
Set-TextValue $ws.Range("D6")  "3.389"
Set-TextValue $ws.Range("D8")  "0.9293"
Set-TextValue $ws.Range("D9")  "0.1428"
Set-TextValue $ws.Range("D10") "0.07434"
Set-TextValue $ws.Range("D11") "0.03377"
Set-TextValue $ws.Range("D12") "0.03038"
Set-TextValue $ws.Range("D13") "0.09354"
Set-TextValue $ws.Range("D14") "3.944"
Set-TextValue $ws.Range("D15") "0.001600"
Set-TextValue $ws.Range("D16") "0.04819"
Set-TextValue $ws.Range("D18") "0.005738"
Set-TextValue $ws.Range("D19") "0.004159"
Set-TextValue $ws.Range("D20") "0.0009852"
Set-TextValue $ws.Range("D23") "6.446"
Set-TextValue $ws.Range("D40") "0.03953"
Set-TextValue $ws.Range("D44") "0.006777"
Set-TextValue $ws.Range("D45") "0.00005202"
Set-TextValue $ws.Range("D49") "0.002274"

# Volume(1h) label (column E) updates (plain text, no conversion risk)
$ws.Range("E41").Value = "40KickTokenKICKBestin24h"
$ws.Range("E44").Value = "43LocalTradersLCT"
